$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 557.8
$ws.Range("I12").Value = 723
$ws.Range("J12").Value = 447.66666
$ws.Range("K12").Value = 723
$ws.Range("L12").Value = 447.66666
$ws.Range("M12").Value = -553
$ws.Range("N12").Value = -787.66666
$ws.Range("H80").Value = 544.1667
$ws.Range("I80").Value = 499.07144
$ws.Range("J80").Value = 607.3
$ws.Range("K80").Value = 1497.21432
$ws.Range("L80").Value = 1821.9
$ws.Range("M80").Value = -499.21432
$ws.Range("N80").Value = -3817.9
$ws.Range("H83").Value = 544.1667
$ws.Range("I83").Value = 499.07144
$ws.Range("J83").Value = 607.3
$ws.Range("K83").Value = 4491.64296
$ws.Range("L83").Value = 5465.7
$ws.Range("M83").Value = 500.3570399999999
$ws.Range("N83").Value = -15449.7
$ws.Range("H94").Value = 16566.334
$ws.Range("I94").Value = 16566.334
$ws.Range("K94").Value = 16566.334
$ws.Range("M94").Value = -16115.334
$ws.Range("H111").Value = 5884119.5
$ws.Range("I111").Value = 2214.125
$ws.Range("J111").Value = 11112480
$ws.Range("K111").Value = 6642.375
$ws.Range("L111").Value = 33337440
$ws.Range("M111").Value = -3575.375
$ws.Range("N111").Value = -33343574
$ws.Range("H112").Value = 2937.2273
$ws.Range("I112").Value = 1363
$ws.Range("J112").Value = 3185.7896
$ws.Range("K112").Value = 4089
$ws.Range("L112").Value = 9557.3688
$ws.Range("M112").Value = -2981
$ws.Range("N112").Value = -11773.3688
$ws.Range("H116").Value = 3024.5
$ws.Range("I116").Value = 2742.2856
$ws.Range("K116").Value = 2742.2856
$ws.Range("M116").Value = 699.7143999999998
$ws.Range("H125").Value = 2625.1667
$ws.Range("I125").Value = 2500.3333
$ws.Range("J125").Value = 2750
$ws.Range("K125").Value = 22502.9997
$ws.Range("L125").Value = 24750
$ws.Range("M125").Value = -20042.9997
$ws.Range("N125").Value = -29670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2018.4166
$ws.Range("I45").Value = 1873.5238
$ws.Range("K45").Value = 1873.5238
$ws.Range("M45").Value = -1496.5238
$ws.Range("H61").Value = 5804.684
$ws.Range("I61").Value = 4731.3335
$ws.Range("J61").Value = 7644.7144
$ws.Range("K61").Value = 4731.3335
$ws.Range("L61").Value = 7644.7144
$ws.Range("M61").Value = -4519.3335
$ws.Range("N61").Value = -8068.7144
$ws.Range("H74").Value = 2459.8948
$ws.Range("I74").Value = 506.94116
$ws.Range("J74").Value = 4040.8572
$ws.Range("K74").Value = 506.94116
$ws.Range("L74").Value = 4040.8572
$ws.Range("M74").Value = 367.05884
$ws.Range("N74").Value = -5788.8572
$ws.Range("H77").Value = 2459.8948
$ws.Range("I77").Value = 506.94116
$ws.Range("J77").Value = 4040.8572
$ws.Range("K77").Value = 2534.7058
$ws.Range("L77").Value = 20204.286
$ws.Range("M77").Value = 1833.2942
$ws.Range("N77").Value = -28940.286
$ws.Range("H97").Value = 1483926.1
$ws.Range("I97").Value = 1545750.2
$ws.Range("K97").Value = 1545750.2
$ws.Range("M97").Value = -1545254.2
$ws.Range("H102").Value = 3027.9524
$ws.Range("I102").Value = 3027.9524
$ws.Range("K102").Value = 3027.9524
$ws.Range("M102").Value = -1405.9524
$ws.Range("H110").Value = 29690092
$ws.Range("I110").Value = 43183730
$ws.Range("J110").Value = 4094.6
$ws.Range("K110").Value = 43183730
$ws.Range("L110").Value = 4094.6
$ws.Range("M110").Value = -43181685
$ws.Range("N110").Value = -8184.6
$ws.Range("H136").Value = 5804.684
$ws.Range("I136").Value = 4731.3335
$ws.Range("J136").Value = 7644.7144
$ws.Range("K136").Value = 14194.0005
$ws.Range("L136").Value = 22934.1432
$ws.Range("M136").Value = -11644.0005
$ws.Range("N136").Value = -28034.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2180.2
$ws.Range("I22").Value = 2662.75
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 2662.75
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = -2489.75
$ws.Range("N22").Value = -596
$ws.Range("H86").Value = 186466.27
$ws.Range("I86").Value = 5136.5557
$ws.Range("K86").Value = 5136.5557
$ws.Range("M86").Value = -4013.5557
$ws.Range("H89").Value = 186466.27
$ws.Range("I89").Value = 5136.5557
$ws.Range("K89").Value = 25682.7785
$ws.Range("M89").Value = -20066.7785
$ws.Range("H99").Value = 2154.889
$ws.Range("I99").Value = 2154.889
$ws.Range("K99").Value = 2154.889
$ws.Range("M99").Value = -656.8890000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34486644
$ws.Range("I31").Value = 111113030
$ws.Range("J31").Value = 4771.3
$ws.Range("K31").Value = 111113030
$ws.Range("L31").Value = 4771.3
$ws.Range("M31").Value = -111112735
$ws.Range("N31").Value = -5361.3
$ws.Range("H34").Value = 34486644
$ws.Range("I34").Value = 111113030
$ws.Range("J34").Value = 4771.3
$ws.Range("K34").Value = 111113030
$ws.Range("L34").Value = 4771.3
$ws.Range("M34").Value = -111112828
$ws.Range("N34").Value = -5175.3
$ws.Range("H58").Value = 3726.76
$ws.Range("I58").Value = 2488.6667
$ws.Range("J58").Value = 5583.9
$ws.Range("K58").Value = 2488.6667
$ws.Range("L58").Value = 5583.9
$ws.Range("M58").Value = -2285.6667
$ws.Range("N58").Value = -5989.9
$ws.Range("H107").Value = 1837.1305
$ws.Range("I107").Value = 1504.5834
$ws.Range("J107").Value = 2199.9092
$ws.Range("K107").Value = 1504.5834
$ws.Range("L107").Value = 2199.9092
$ws.Range("M107").Value = 415.4166
$ws.Range("N107").Value = -6039.9092
$ws.Range("H132").Value = 84602.53
$ws.Range("I132").Value = 3681.6667
$ws.Range("K132").Value = 11045.0001
$ws.Range("M132").Value = -8515.000100000001
$ws.Range("H136").Value = 3726.76
$ws.Range("I136").Value = 2488.6667
$ws.Range("J136").Value = 5583.9
$ws.Range("K136").Value = 7466.000100000001
$ws.Range("L136").Value = 16751.7
$ws.Range("M136").Value = -4916.000100000001
$ws.Range("N136").Value = -21851.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18500.092
$ws.Range("I70").Value = 15944.889
$ws.Range("K70").Value = 15944.889
$ws.Range("M70").Value = -15674.889
$ws.Range("H73").Value = 18500.092
$ws.Range("I73").Value = 15944.889
$ws.Range("K73").Value = 15944.889
$ws.Range("M73").Value = -15008.889
$ws.Range("H113").Value = 22046.295
$ws.Range("I113").Value = 1756.1111
$ws.Range("K113").Value = 1756.1111
$ws.Range("M113").Value = 413.8888999999999
$ws.Range("H126").Value = 2568.875
$ws.Range("J126").Value = 3409.3333
$ws.Range("L126").Value = 10227.9999
$ws.Range("N126").Value = -15167.9999
$ws.Range("H132").Value = 3888.7104
$ws.Range("I132").Value = 1916.3182
$ws.Range("J132").Value = 6600.75
$ws.Range("K132").Value = 5748.9546
$ws.Range("L132").Value = 19802.25
$ws.Range("M132").Value = -3218.9546
$ws.Range("N132").Value = -24862.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2769.8718
$ws.Range("I22").Value = 1823.3334
$ws.Range("J22").Value = 4284.3335
$ws.Range("K22").Value = 1823.3334
$ws.Range("L22").Value = 4284.3335
$ws.Range("M22").Value = -1528.3334
$ws.Range("N22").Value = -4874.3335
$ws.Range("H27").Value = 2769.8718
$ws.Range("I27").Value = 1823.3334
$ws.Range("J27").Value = 4284.3335
$ws.Range("K27").Value = 1823.3334
$ws.Range("L27").Value = 4284.3335
$ws.Range("M27").Value = -1716.3334
$ws.Range("N27").Value = -4498.3335
$ws.Range("H46").Value = 3850.7646
$ws.Range("I46").Value = 1741.6666
$ws.Range("J46").Value = 4302.7144
$ws.Range("K46").Value = 1741.6666
$ws.Range("L46").Value = 4302.7144
$ws.Range("M46").Value = -1553.6666
$ws.Range("N46").Value = -4678.7144
$ws.Range("H55").Value = 331.7647
$ws.Range("I55").Value = 282.64285
$ws.Range("K55").Value = 282.64285
$ws.Range("M55").Value = -109.64285
$ws.Range("H61").Value = 2171.5557
$ws.Range("I61").Value = 1377
$ws.Range("J61").Value = 4952.5
$ws.Range("K61").Value = 1377
$ws.Range("L61").Value = 4952.5
$ws.Range("M61").Value = -1175
$ws.Range("N61").Value = -5356.5
$ws.Range("H82").Value = 2497.125
$ws.Range("I82").Value = 1970
$ws.Range("J82").Value = 3024.25
$ws.Range("K82").Value = 1970
$ws.Range("L82").Value = 3024.25
$ws.Range("M82").Value = -1609
$ws.Range("N82").Value = -3746.25
$ws.Range("H85").Value = 2497.125
$ws.Range("I85").Value = 1970
$ws.Range("J85").Value = 3024.25
$ws.Range("K85").Value = 1970
$ws.Range("L85").Value = 3024.25
$ws.Range("M85").Value = -722
$ws.Range("N85").Value = -5520.25
$ws.Range("H93").Value = 5801.7144
$ws.Range("J93").Value = 6495.8335
$ws.Range("L93").Value = 6495.8335
$ws.Range("N93").Value = -8991.833500000001
$ws.Range("H113").Value = 2171.5557
$ws.Range("I113").Value = 1377
$ws.Range("J113").Value = 4952.5
$ws.Range("K113").Value = 1377
$ws.Range("L113").Value = 4952.5
$ws.Range("M113").Value = 793
$ws.Range("N113").Value = -9292.5
$ws.Range("H132").Value = 5150.154
$ws.Range("I132").Value = 4278.722
$ws.Range("J132").Value = 7110.875
$ws.Range("K132").Value = 12836.166
$ws.Range("L132").Value = 21332.625
$ws.Range("M132").Value = -10306.166
$ws.Range("N132").Value = -26392.625
$ws.Range("H136").Value = 3973.311
$ws.Range("I136").Value = 3105.524
$ws.Range("J136").Value = 4732.625
$ws.Range("K136").Value = 9316.572
$ws.Range("L136").Value = 14197.875
$ws.Range("M136").Value = -6766.572
$ws.Range("N136").Value = -19297.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4778.5757
$ws.Range("I132").Value = 2525.3333
$ws.Range("J132").Value = 7482.467
$ws.Range("K132").Value = 7575.999899999999
$ws.Range("L132").Value = 22447.401
$ws.Range("M132").Value = -5045.999899999999
$ws.Range("N132").Value = -27507.401
$ws.Range("H136").Value = 6327
$ws.Range("I136").Value = 2683
$ws.Range("J136").Value = 10699.8
$ws.Range("K136").Value = 8049
$ws.Range("L136").Value = 32099.4
$ws.Range("M136").Value = -5499
$ws.Range("N136").Value = -37199.39999999999
